$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (active) on row 2 currently holds the numeric value 1 (style 2,
# numFmtId 1 = "0"). Replace it with the literal text "true" while keeping
# the same cell style. A leading apostrophe forces Excel to store it as a
# text string (new shared string) rather than re-coercing it back into a
# number or auto-detecting a Boolean.
$ws.Range("C2").Value = "'true"

# Move the sheet's active selection from E3 to A3.
$ws.Range("A3").Select() | Out-Null
